$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2467.4211
$ws.Range("I28").Value = 491.7143
$ws.Range("K28").Value = 491.7143
$ws.Range("M28").Value = -6.71429999999998
$ws.Range("H40").Value = 14535.454
$ws.Range("I40").Value = 23800.8
$ws.Range("J40").Value = 6814.3335
$ws.Range("K40").Value = 23800.8
$ws.Range("L40").Value = 6814.3335
$ws.Range("M40").Value = -23625.8
$ws.Range("N40").Value = -7164.3335
$ws.Range("H76").Value = 6759.6
$ws.Range("I76").Value = 5519.2
$ws.Range("K76").Value = 5519.2
$ws.Range("M76").Value = -5204.2
$ws.Range("H79").Value = 6759.6
$ws.Range("I79").Value = 5519.2
$ws.Range("K79").Value = 5519.2
$ws.Range("M79").Value = -4427.2
$ws.Range("H98").Value = 10831.667
$ws.Range("I98").Value = 10831.667
$ws.Range("K98").Value = 10831.667
$ws.Range("M98").Value = -9333.666999999999
$ws.Range("H122").Value = 10831.667
$ws.Range("I122").Value = 10831.667
$ws.Range("K122").Value = 32495.001
$ws.Range("M122").Value = -30045.001
$ws.Range("H125").Value = 9011676
$ws.Range("J125").Value = 10419516
$ws.Range("L125").Value = 93775644
$ws.Range("N125").Value = -93780564
$ws.Range("H127").Value = 1201.6666
$ws.Range("I127").Value = 856.36365
$ws.Range("K127").Value = 2569.09095
$ws.Range("M127").Value = 2390.90905
$ws.Range("H132").Value = 26318836
$ws.Range("I132").Value = 28574538
$ws.Range("J132").Value = 2316
$ws.Range("K132").Value = 85723614
$ws.Range("L132").Value = 6948
$ws.Range("M132").Value = -85721084
$ws.Range("N132").Value = -12008
$ws.Range("H138").Value = 4041.8572
$ws.Range("I138").Value = 2524.375
$ws.Range("J138").Value = 4491.4814
$ws.Range("K138").Value = 7573.125
$ws.Range("L138").Value = 13474.4442
$ws.Range("M138").Value = -2433.125
$ws.Range("N138").Value = -23754.4442

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 2263.7273
$ws.Range("I28").Value = 2263.7273
$ws.Range("K28").Value = 2263.7273
$ws.Range("M28").Value = -2071.7273
$ws.Range("H32").Value = 9801.808000000001
$ws.Range("I32").Value = 8329.944
$ws.Range("K32").Value = 8329.944
$ws.Range("M32").Value = -8042.944
$ws.Range("H45").Value = 4798108.5
$ws.Range("I45").Value = 6255968.5
$ws.Range("J45").Value = 7998.7144
$ws.Range("K45").Value = 6255968.5
$ws.Range("L45").Value = 7998.7144
$ws.Range("M45").Value = -6255591.5
$ws.Range("N45").Value = -8752.714400000001
$ws.Range("H61").Value = 11915.167
$ws.Range("I61").Value = 12589.272
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 12589.272
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -12377.272
$ws.Range("N61").Value = -4924
$ws.Range("H63").Value = 5818.5
$ws.Range("I63").Value = 1230.75
$ws.Range("K63").Value = 1230.75
$ws.Range("M63").Value = -544.75
$ws.Range("H66").Value = 5818.5
$ws.Range("I66").Value = 1230.75
$ws.Range("K66").Value = 6153.75
$ws.Range("M66").Value = -2721.75
$ws.Range("H99").Value = 2263.7273
$ws.Range("I99").Value = 2263.7273
$ws.Range("K99").Value = 2263.7273
$ws.Range("M99").Value = 731.2727
$ws.Range("H110").Value = 1737611.1
$ws.Range("I110").Value = 5558039.5
$ws.Range("J110").Value = 1052.7273
$ws.Range("K110").Value = 5558039.5
$ws.Range("L110").Value = 1052.7273
$ws.Range("M110").Value = -5555994.5
$ws.Range("N110").Value = -5142.7273
$ws.Range("H122").Value = 550370.2
$ws.Range("I122").Value = 1650.6897
$ws.Range("K122").Value = 4952.0691
$ws.Range("M122").Value = -2502.0691
$ws.Range("H130").Value = 49051.332
$ws.Range("J130").Value = 49051.332
$ws.Range("L130").Value = 49051.332
$ws.Range("N130").Value = -59091.332
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 7664
$ws.Range("I132").Value = 8487.058999999999
$ws.Range("J132").Value = 6109.3335
$ws.Range("K132").Value = 25461.177
$ws.Range("L132").Value = 18328.0005
$ws.Range("M132").Value = -22931.177
$ws.Range("N132").Value = -23388.0005
$ws.Range("H136").Value = 11915.167
$ws.Range("I136").Value = 12589.272
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 37767.81600000001
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -35217.81600000001
$ws.Range("N136").Value = -18600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1044.3077
$ws.Range("J22").Value = 1217
$ws.Range("L22").Value = 1217
$ws.Range("N22").Value = -1563
$ws.Range("H88").Value = 29893.143
$ws.Range("J88").Value = 29893.143
$ws.Range("L88").Value = 29893.143
$ws.Range("N88").Value = -30705.143
$ws.Range("H91").Value = 29893.143
$ws.Range("J91").Value = 29893.143
$ws.Range("L91").Value = 29893.143
$ws.Range("N91").Value = -32701.143
$ws.Range("H94").Value = 2273551.5
$ws.Range("J94").Value = 1563.7778
$ws.Range("L94").Value = 1563.7778
$ws.Range("N94").Value = -2465.7778

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 257.04544
$ws.Range("I7").Value = 129.8
$ws.Range("J7").Value = 529.7143
$ws.Range("K7").Value = 129.8
$ws.Range("L7").Value = 529.7143
$ws.Range("M7").Value = -16.80000000000001
$ws.Range("N7").Value = -755.7143
$ws.Range("H31").Value = 15091.685
$ws.Range("I31").Value = 7787
$ws.Range("K31").Value = 7787
$ws.Range("M31").Value = -7492
$ws.Range("H34").Value = 15091.685
$ws.Range("I34").Value = 7787
$ws.Range("K34").Value = 7787
$ws.Range("M34").Value = -7585
$ws.Range("H97").Value = 24500
$ws.Range("J97").Value = 24500
$ws.Range("L97").Value = 24500
$ws.Range("N97").Value = -26482
$ws.Range("H107").Value = 1178.7368
$ws.Range("I107").Value = 1119.2858
$ws.Range("K107").Value = 1119.2858
$ws.Range("M107").Value = 800.7141999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 89061
$ws.Range("J12").Value = 183.66667
$ws.Range("L12").Value = 551.00001
$ws.Range("N12").Value = -897.00001
$ws.Range("H37").Value = 66799.14
$ws.Range("J37").Value = 66799.14
$ws.Range("L37").Value = 200397.42
$ws.Range("N37").Value = -200621.42
$ws.Range("H131").Value = 13891036
$ws.Range("J131").Value = 12348027
$ws.Range("L131").Value = 37044081
$ws.Range("N131").Value = -37054161
$ws.Range("H136").Value = 1733.25
$ws.Range("I136").Value = 1144.3334
$ws.Range("K136").Value = 3433.0002
$ws.Range("M136").Value = 1666.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 100003750
$ws.Range("I70").Value = 200000000
$ws.Range("J70").Value = 7500
$ws.Range("K70").Value = 200000000
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -199999730
$ws.Range("N70").Value = -8040
$ws.Range("H73").Value = 100003750
$ws.Range("I73").Value = 200000000
$ws.Range("J73").Value = 7500
$ws.Range("K73").Value = 200000000
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -199999064
$ws.Range("N73").Value = -9372
$ws.Range("H102").Value = 4689957.5
$ws.Range("I102").Value = 7938316.5
$ws.Range("J102").Value = 1658155.6
$ws.Range("K102").Value = 7938316.5
$ws.Range("L102").Value = 1658155.6
$ws.Range("M102").Value = -7936694.5
$ws.Range("N102").Value = -1661399.6
$ws.Range("H122").Value = 472398.62
$ws.Range("I122").Value = 560229.4399999999
$ws.Range("J122").Value = 3967.6667
$ws.Range("K122").Value = 1680688.32
$ws.Range("L122").Value = 11903.0001
$ws.Range("M122").Value = -1678238.32
$ws.Range("N122").Value = -16803.0001
$ws.Range("H126").Value = 3791025
$ws.Range("I126").Value = 2068908.9
$ws.Range("J126").Value = 6948238
$ws.Range("K126").Value = 6206726.699999999
$ws.Range("L126").Value = 20844714
$ws.Range("M126").Value = -6204256.699999999
$ws.Range("N126").Value = -20849654
$ws.Range("H132").Value = 9887.645500000001
$ws.Range("I132").Value = 7263.2607
$ws.Range("K132").Value = 21789.7821
$ws.Range("M132").Value = -19259.7821

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 15152997
$ws.Range("I93").Value = 18519812
$ws.Range("K93").Value = 18519812
$ws.Range("M93").Value = -18518564
$ws.Range("H132").Value = 22275.1
$ws.Range("I132").Value = 26781.5
$ws.Range("K132").Value = 80344.5
$ws.Range("M132").Value = -77814.5
$ws.Range("H140").Value = 200429
$ws.Range("J140").Value = 200429
$ws.Range("L140").Value = 200429
$ws.Range("N140").Value = -210789

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2194.9524
$ws.Range("I122").Value = 1718.6875
$ws.Range("J122").Value = 3719
$ws.Range("K122").Value = 5156.0625
$ws.Range("L122").Value = 11157
$ws.Range("M122").Value = -2706.0625
$ws.Range("N122").Value = -16057
$ws.Range("H126").Value = 3060.3845
$ws.Range("I126").Value = 2988.7
$ws.Range("K126").Value = 8966.099999999999
$ws.Range("M126").Value = -6496.099999999999
$ws.Range("H136").Value = 5553.7295
$ws.Range("J136").Value = 3696.2632
$ws.Range("L136").Value = 11088.7896
$ws.Range("N136").Value = -16188.7896
